# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.102.32"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "2.962.87"
$ws.Range("E3").Value = "  +0.80%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "380.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("E7").Value = "  +1.84%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.10%  "

$ws.Range("E11").Value = "  -0.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0853"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.14%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.82%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.426.90"
$ws.Range("E14").Value = "  +0.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "12.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +73.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.69%  "

$ws.Range("D17").Value = "2.965.66"
$ws.Range("E17").Value = "  +0.98%  "

$ws.Range("E18").Value = "  +3.89%  "

$ws.Range("D19").Value = "51.163.72"
$ws.Range("E19").Value = "  +0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.93%  "

$ws.Range("D22").Value = "0.0₃0959"
$ws.Range("E22").Value = "  +1.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +16.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.109"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "51.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "34.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.71%  "

$ws.Range("E35").Value = "  +2.19%  "

$ws.Range("E36").Value = "  -4.15%  "

$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.72%  "

$ws.Range("E39").Value = "  +2.27%  "

$ws.Range("E40").Value = "  +1.60%  "

$ws.Range("E41").Value = "  +3.06%  "

$ws.Range("E42").Value = "  -2.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.84%  "

$ws.Range("D46").Value = "2.091.83"
$ws.Range("E46").Value = "  +4.23%  "

$ws.Range("E47").Value = "  -0.89%  "

$ws.Range("E48").Value = "  -0.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.261"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0322"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.99%  "

